$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new bonus site ("nimacasino") needs to be inserted in its correct
# alphabetically-sorted position. In the existing sheet that slot is
# row 324 (currently "nisanbet"), so push everything from row 324 down
# by one row and populate the freed row with the new entry.
$ws.Rows.Item(324).Insert()

$ws.Cells.Item(324, 1).Value = "nimacasino"
$ws.Cells.Item(324, 2).Value = "Maks 1k çekim"
$ws.Cells.Item(324, 3).Value = "yatırımsız"
$ws.Cells.Item(324, 4).Value = "Evet"

# Restore the view to the top of the sheet/table.
$sheetView = $ws.Application.ActiveWindow
$sheetView.ScrollRow = 1
$ws.Range("A2").Select()
